$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.979.13"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "1.786.26"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'" + "221.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'" + "31.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.41%  "

$ws.Range("E9").Value = "  +1.38%  "

$ws.Range("D10").Value = "'" + "0.0713"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.51%  "

$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "2.041.02"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").Value = "1.784.23"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "'" + "10.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.38%  "

$ws.Range("D15").Value = "'" + "0.628"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "33.957.29"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "'" + "4.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "'" + "67.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").Value = "'" + "244.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").Value = "'" + "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "'" + "10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.85%  "

$ws.Range("D23").Value = "'" + "4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("D25").Value = "'" + "157.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("D26").Value = "'" + "16.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").Value = "'" + "7.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  -1.79%  "

$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").Value = "1.410.21"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").Value = "'" + "0.642"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("D39").Value = "'" + "0.938"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.20%  "

$ws.Range("D40").Value = "'" + "79.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.05%  "

$ws.Range("E41").Value = "  -3.30%  "

$ws.Range("E42").Value = "  -0.95%  "

$ws.Range("E43").Value = "  +2.08%  "

$ws.Range("D44").Value = "'" + "5.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").Value = "'" + "0.0493"
$ws.Range("D45").Style = "Normal"

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.940.01"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'" + "1.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.00%  "

$ws.Range("D48").Value = "'" + "105.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").Value = "'" + "11.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("E51").Value = "  -1.01%  "
